$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("monster")

# --- Column J : c4 / array<int> ---
$ws.Range("J1").Value = "c4"
$ws.Range("J2").Value = "array<int>"
$ws.Range("J4").Value = "11,21"
$ws.Range("J5").Value = "11,11,11,2"
$ws.Range("J6").Value = 3
$ws.Range("J7").Value = 4
$ws.Range("J8").Value = 1
$ws.Range("J9").Value = 1
$ws.Range("J10").Value = 1
$ws.Range("J11").Value = 1
$ws.Range("J12").Value = 1
$ws.Range("J13").Value = 1

# --- Column K : c5 / group<int> ---
$ws.Range("K1").Value = "c5"
$ws.Range("K2").Value = "group<int>"
$ws.Range("K4").Value = "[11,12,13],[11,12,14]"
$ws.Range("K5").Value = "[11,12,13],[11,12,15]"
$ws.Range("K6").Value = "[11,12,13],[11,12,16]"
$ws.Range("K7").Value = "[11,12,13],[11,12,17]"
$ws.Range("K8").Value = "[11,12,13],[11,12,18]"
$ws.Range("K9").Value = "[11,12,13],[11,12,19]"
$ws.Range("K10").Value = "[11,12,13],[11,12,20]"
$ws.Range("K11").Value = "[11,12,13],[11,12,21]"
$ws.Range("K12").Value = "[11,12,13],[11,12,22]"
$ws.Range("K13").Value = "[11,12,13],[11,12,23]"

# --- Column L : c8 / array<bool> (header L1 written last, see bottom) ---
$ws.Range("L2").Value = "array<bool>"
$ws.Range("L4").Value = "[true,false,true]"
$ws.Range("L5").Value = "[true,true]"
$ws.Range("L6").Value = "[true,false,true]"
$ws.Range("L7").Value = "[true,false,true]"
$ws.Range("L8").Value = "[true,false,true]"
$ws.Range("L9").Value = "[true,false,true]"
$ws.Range("L10").Value = "[true]"
$ws.Range("L11").Value = "[true,false,true]"
$ws.Range("L12").Value = "[true,false,true]"
$ws.Range("L13").Value = "[true,false,true]"

# --- Column M : c6 / bool ---
$ws.Range("M1").Value = "c6"
$ws.Range("M2").Value = "bool"
$ws.Range("M4").Value = $true
$ws.Range("M5").Value = $false
$ws.Range("M6").Value = 1
$ws.Range("M7").Value = 2
$ws.Range("M8").Value = 4
$ws.Range("M9").Value = 0
$ws.Range("M10").Value = 0
$ws.Range("M11").Value = 0
$ws.Range("M12").Value = 0
$ws.Range("M13").Value = 0

# --- Column N : c7 / group<bool> ---
$ws.Range("N1").Value = "c7"
$ws.Range("N2").Value = "group<bool>"
$ws.Range("N4").Value = "[true,false],[true,false],[true,false],[true,false]"
$ws.Range("N5").Value = "[true,false],[true,false],[true,false],[true,false]"
$ws.Range("N6").Value = "[true,false],[true,false],[true,false],[true,false]"
$ws.Range("N7").Value = "[true,false],[true,false],[true,false],[true,false]"
$ws.Range("N8").Value = "[true,false],[true,false],[true,false],[true,false]"
$ws.Range("N9").Value = "[true,false],[true,false],[true,false],[true,false]"
$ws.Range("N10").Value = "[true,false],[true,false],[true,false],[true,false]"
$ws.Range("N11").Value = "[true,false],[true,false],[true,false],[true,false]"
$ws.Range("N12").Value = "[true,false],[true,false],[true,false],[true,false]"
$ws.Range("N13").Value = "[true,false],[true,false],[true,false],[true,false]"

# L1 header is added last (matches authoring order captured in the shared-string table)
$ws.Range("L1").Value = "c8"

# --- Header row formatting: new header cells (J1:N1) match the existing
#     header row's vertical-center style used by A1:I1 ---
$ws.Range("J1:N1").VerticalAlignment = -4108

# --- Column widths (character units); new/changed columns A, B, C, I, K, L, N ---
$ws.Columns.Item(1).ColumnWidth = 9.857142857142858
$ws.Columns.Item(2).ColumnWidth = 12.428571428571429
$ws.Columns.Item(3).ColumnWidth = 14.0
$ws.Columns.Item(9).ColumnWidth = 29.857142857142858
$ws.Columns.Item(11).ColumnWidth = 31.857142857142858
$ws.Columns.Item(12).ColumnWidth = 16.0
$ws.Columns.Item(14).ColumnWidth = 52.714285714285715

# --- View state: scroll / selection as left by the editing session ---
$ws.Range("L23").Select()
